{"js": "const pairs = [\n  [\"2025-02-20 Thursday\", \"2025-02-21 Friday\"],\n  [\"33\u00d772=\", \"35\u00d779=\"],\n  [\"72\u00d765=\", \"45\u00d798=\"],\n  [\"55\u00d786=\", \"43\u00d732=\"],\n  [\"45\u00d767=\", \"14\u00d799=\"],\n  [\"92\u00d712=\", \"90\u00d752=\"],\n  [\"81\u00d731=\", \"92\u00d740=\"],\n  [\"79\u00d750=\", \"94\u00d797=\"],\n  [\"36\u00d793=\", \"33\u00d719=\"],\n  [\"41\u00d762=\", \"11\u00d792=\"],\n  [\"62\u00d716=\", \"30\u00d773=\"],\n  [\"37\u00d789=\", \"27\u00d745=\"],\n  [\"70\u00d712=\", \"64\u00d755=\"],\n  [\"84\u00d774=\", \"12\u00d773=\"],\n  [\"64\u00d792=\", \"69\u00d724=\"],\n  [\"18\u00d761=\", \"41\u00d790=\"],\n  [\"29\u00d776=\", \"38\u00d733=\"],\n  [\"17\u00d744=\", \"34\u00d738=\"],\n  [\"42\u00d715=\", \"84\u00d743=\"],\n  [\"75\u00d753=\", \"66\u00d796=\"],\n  [\"73\u00d732=\", \"59\u00d769=\"],\n  [\"83\u00d720=\", \"98\u00d750=\"],\n  [\"37\u00d718=\", \"53\u00d740=\"],\n  [\"48\u00d770=\", \"36\u00d767=\"],\n  [\"58\u00d776=\", \"34\u00d799=\"],\n  [\"77\u00d712=\", \"53\u00d730=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-02-20 Thursday\", \"2025-02-21 Friday\"),\n    @(\"33\u00d772=\", \"35\u00d779=\"),\n    @(\"72\u00d765=\", \"45\u00d798=\"),\n    @(\"55\u00d786=\", \"43\u00d732=\"),\n    @(\"45\u00d767=\", \"14\u00d799=\"),\n    @(\"92\u00d712=\", \"90\u00d752=\"),\n    @(\"81\u00d731=\", \"92\u00d740=\"),\n    @(\"79\u00d750=\", \"94\u00d797=\"),\n    @(\"36\u00d793=\", \"33\u00d719=\"),\n    @(\"41\u00d762=\", \"11\u00d792=\"),\n    @(\"62\u00d716=\", \"30\u00d773=\"),\n    @(\"37\u00d789=\", \"27\u00d745=\"),\n    @(\"70\u00d712=\", \"64\u00d755=\"),\n    @(\"84\u00d774=\", \"12\u00d773=\"),\n    @(\"64\u00d792=\", \"69\u00d724=\"),\n    @(\"18\u00d761=\", \"41\u00d790=\"),\n    @(\"29\u00d776=\", \"38\u00d733=\"),\n    @(\"17\u00d744=\", \"34\u00d738=\"),\n    @(\"42\u00d715=\", \"84\u00d743=\"),\n    @(\"75\u00d753=\", \"66\u00d796=\"),\n    @(\"73\u00d732=\", \"59\u00d769=\"),\n    @(\"83\u00d720=\", \"98\u00d750=\"),\n    @(\"37\u00d718=\", \"53\u00d740=\"),\n    @(\"48\u00d770=\", \"36\u00d767=\"),\n    @(\"58\u00d776=\", \"34\u00d799=\"),\n    @(\"77\u00d712=\", \"53\u00d730=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)  # wdReplaceAll\n}\n"}
